$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 45033
$ws.Range("L2").Value = 'Especial'
$ws.Range("N2").Value = 13000
$ws.Range("O2").Value = 13000
$ws.Range("P2").Value = 13000
$ws.Range("R2").Value = 'Región de O''Higgins'
$ws.Range("S2").Value = 722

$ws.Range("D3").Value = 45033
$ws.Range("L3").Value = 'Primera'
$ws.Range("M3").Value = 80
$ws.Range("N3").Value = 12000
$ws.Range("O3").Value = 12000
$ws.Range("P3").Value = 12000
$ws.Range("S3").Value = 667

$ws.Range("D4").Value = 45070
$ws.Range("M4").Value = 60
$ws.Range("N4").Value = 10000
$ws.Range("O4").Value = 10000
$ws.Range("P4").Value = 10000
$ws.Range("S4").Value = 556

$ws.Range("D5").Value = 44699
$ws.Range("M5").Value = 60
$ws.Range("Q5").Value = '$/caja 15 kilos granel'
$ws.Range("R5").Value = 'Provincia de Curicó'
$ws.Range("S5").Value = 867
$ws.Range("T5").Value = 15

$ws.Range("D6").Value = 44699
$ws.Range("M6").Value = 120
$ws.Range("N6").Value = 11000
$ws.Range("P6").Value = 11500
$ws.Range("Q6").Value = '$/caja 15 kilos granel'
$ws.Range("R6").Value = 'Provincia de Curicó'
$ws.Range("S6").Value = 767
$ws.Range("T6").Value = 15

$ws.Range("D7").Value = 45089
$ws.Range("N7").Value = 11000
$ws.Range("O7").Value = 11000
$ws.Range("P7").Value = 11000
$ws.Range("Q7").Value = '$/caja 18 kilos empedrada'
$ws.Range("R7").Value = 'Región del Maule'
$ws.Range("S7").Value = 611
$ws.Range("T7").Value = 18

$ws.Range("D8").Value = 45089
$ws.Range("M8").Value = 50
$ws.Range("N8").Value = 9000
$ws.Range("O8").Value = 9000
$ws.Range("P8").Value = 9000
$ws.Range("Q8").Value = '$/caja 18 kilos empedrada'
$ws.Range("R8").Value = 'Región del Maule'
$ws.Range("S8").Value = 500
$ws.Range("T8").Value = 18

$ws.Range("D9").Value = 45089
$ws.Range("L9").Value = 'Segunda'
$ws.Range("M9").Value = 30
$ws.Range("N9").Value = 7000
$ws.Range("O9").Value = 7000
$ws.Range("P9").Value = 7000
$ws.Range("R9").Value = 'Región del Maule'
$ws.Range("S9").Value = 389

$ws.Range("D10").Value = 45040
$ws.Range("L10").Value = 'Especial'
$ws.Range("M10").Value = 50
$ws.Range("N10").Value = 13000
$ws.Range("O10").Value = 13000
$ws.Range("P10").Value = 13000
$ws.Range("S10").Value = 722

$ws.Range("D11").Value = 45040

$ws.Range("D12").Value = 45069
$ws.Range("L12").Value = 'Primera'
$ws.Range("M12").Value = 60
$ws.Range("N12").Value = 12000
$ws.Range("O12").Value = 12000
$ws.Range("P12").Value = 12000
$ws.Range("S12").Value = 667

$ws.Range("D13").Value = 45069
$ws.Range("L13").Value = 'Segunda'
$ws.Range("M13").Value = 40
$ws.Range("N13").Value = 10000
$ws.Range("O13").Value = 10000
$ws.Range("P13").Value = 10000
$ws.Range("S13").Value = 556

$ws.Range("D14").Value = 45062
$ws.Range("L14").Value = 'Especial'
$ws.Range("M14").Value = 50
$ws.Range("N14").Value = 13000
$ws.Range("O14").Value = 13000
$ws.Range("P14").Value = 13000
$ws.Range("S14").Value = 722

$ws.Range("D15").Value = 45062
$ws.Range("L15").Value = 'Primera'

$ws.Range("D16").Value = 45071
$ws.Range("M16").Value = 40
$ws.Range("N16").Value = 12000
$ws.Range("O16").Value = 12000
$ws.Range("P16").Value = 12000
$ws.Range("S16").Value = 667

$ws.Range("D17").Value = 45071
$ws.Range("M17").Value = 40
$ws.Range("N17").Value = 10000
$ws.Range("O17").Value = 10000
$ws.Range("P17").Value = 10000
$ws.Range("S17").Value = 556

$ws.Range("D18").Value = 45049
$ws.Range("L18").Value = 'Especial'
$ws.Range("M18").Value = 50
$ws.Range("N18").Value = 13000
$ws.Range("O18").Value = 13000
$ws.Range("P18").Value = 13000
$ws.Range("S18").Value = 722

$ws.Range("D19").Value = 45049
$ws.Range("L19").Value = 'Primera'
$ws.Range("M19").Value = 60
$ws.Range("N19").Value = 12000
$ws.Range("O19").Value = 12000
$ws.Range("P19").Value = 12000
$ws.Range("S19").Value = 667

$ws.Range("D20").Value = 45020
$ws.Range("L20").Value = 'Primera'
$ws.Range("M20").Value = 60
$ws.Range("N20").Value = 12000
$ws.Range("O20").Value = 12000
$ws.Range("P20").Value = 12000
$ws.Range("Q20").Value = '$/caja 18 kilos granel'
$ws.Range("R20").Value = 'Región de O''Higgins'
$ws.Range("S20").Value = 667

$ws.Range("D21").Value = 45076
$ws.Range("M21").Value = 30
$ws.Range("N21").Value = 12000
$ws.Range("O21").Value = 12000
$ws.Range("P21").Value = 12000
$ws.Range("Q21").Value = '$/caja 15 kilos granel'
$ws.Range("R21").Value = 'Región de O''Higgins'
$ws.Range("S21").Value = 800
$ws.Range("T21").Value = 15

$ws.Range("D22").Value = 45076
$ws.Range("M22").Value = 30
$ws.Range("N22").Value = 10000
$ws.Range("O22").Value = 10000
$ws.Range("P22").Value = 10000
$ws.Range("Q22").Value = '$/caja 15 kilos granel'
$ws.Range("R22").Value = 'Región de O''Higgins'
$ws.Range("S22").Value = 667
$ws.Range("T22").Value = 15

$ws.Range("D23").Value = 45050
$ws.Range("L23").Value = 'Especial'
$ws.Range("M23").Value = 50
$ws.Range("N23").Value = 13000
$ws.Range("O23").Value = 13000
$ws.Range("P23").Value = 13000
$ws.Range("Q23").Value = '$/caja 18 kilos empedrada'
$ws.Range("S23").Value = 722

$ws.Range("D24").Value = 45050
$ws.Range("L24").Value = 'Primera'
$ws.Range("M24").Value = 40
$ws.Range("N24").Value = 12000
$ws.Range("O24").Value = 12000
$ws.Range("P24").Value = 12000
$ws.Range("S24").Value = 667

$ws.Range("D25").Value = 45044
$ws.Range("L25").Value = 'Especial'
$ws.Range("N25").Value = 13000
$ws.Range("O25").Value = 13000
$ws.Range("P25").Value = 13000
$ws.Range("S25").Value = 722

$ws.Range("D26").Value = 45044
$ws.Range("L26").Value = 'Primera'
$ws.Range("M26").Value = 40
$ws.Range("N26").Value = 12000
$ws.Range("O26").Value = 12000
$ws.Range("P26").Value = 12000
$ws.Range("R26").Value = 'Región de O''Higgins'
$ws.Range("S26").Value = 667

$ws.Range("D27").Value = 45079
$ws.Range("L27").Value = 'Especial'
$ws.Range("N27").Value = 12000
$ws.Range("O27").Value = 12000
$ws.Range("P27").Value = 12000
$ws.Range("R27").Value = 'Región de O''Higgins'
$ws.Range("S27").Value = 667

$ws.Range("D28").Value = 45079
$ws.Range("L28").Value = 'Primera'
$ws.Range("N28").Value = 10000
$ws.Range("O28").Value = 10000
$ws.Range("P28").Value = 10000
$ws.Range("R28").Value = 'Región de O''Higgins'
$ws.Range("S28").Value = 556

$ws.Range("D29").Value = 45079
$ws.Range("L29").Value = 'Segunda'
$ws.Range("M29").Value = 20
$ws.Range("N29").Value = 9000
$ws.Range("O29").Value = 9000
$ws.Range("P29").Value = 9000
$ws.Range("S29").Value = 500

$ws.Range("D30").Value = 45043
$ws.Range("M30").Value = 40

$ws.Range("D31").Value = 45043

$ws.Range("D32").Value = 45090
$ws.Range("L32").Value = 'Especial'
$ws.Range("M32").Value = 80
$ws.Range("N32").Value = 11000
$ws.Range("O32").Value = 11000
$ws.Range("P32").Value = 11000
$ws.Range("S32").Value = 611

$ws.Range("D33").Value = 45090
$ws.Range("L33").Value = 'Primera'
$ws.Range("M33").Value = 50
$ws.Range("N33").Value = 10000
$ws.Range("O33").Value = 10000
$ws.Range("P33").Value = 10000
$ws.Range("R33").Value = 'Región del Maule'
$ws.Range("S33").Value = 556

$ws.Range("D34").Value = 45090
$ws.Range("L34").Value = 'Segunda'
$ws.Range("M34").Value = 40
$ws.Range("N34").Value = 8000
$ws.Range("O34").Value = 8000
$ws.Range("P34").Value = 8000
$ws.Range("R34").Value = 'Región del Maule'
$ws.Range("S34").Value = 444

$ws.Range("D35").Value = 45085
$ws.Range("M35").Value = 50
$ws.Range("N35").Value = 10000
$ws.Range("O35").Value = 10000
$ws.Range("P35").Value = 10000
$ws.Range("Q35").Value = '$/caja 18 kilos empedrada'
$ws.Range("R35").Value = 'Región del Maule'
$ws.Range("S35").Value = 556
$ws.Range("T35").Value = 18

$ws.Range("D36").Value = 45021
$ws.Range("L36").Value = 'Primera'
$ws.Range("M36").Value = 50
$ws.Range("N36").Value = 12000
$ws.Range("O36").Value = 12000
$ws.Range("P36").Value = 12000
$ws.Range("Q36").Value = '$/caja 18 kilos granel'
$ws.Range("T36").Value = 18

$ws.Range("D37").Value = 45083
$ws.Range("M37").Value = 60
$ws.Range("N37").Value = 9000
$ws.Range("O37").Value = 10000
$ws.Range("P37").Value = 9500
$ws.Range("Q37").Value = '$/caja 18 kilos empedrada'
$ws.Range("R37").Value = 'Región del Maule'
$ws.Range("S37").Value = 528
